# Apply updated crypto price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.212.75'
$ws.Range('E2').Value = '  +3.02%  '
$ws.Range('D3').Value = '1.896.06'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = "'324.96"
$ws.Range('E5').Value = '  +3.28%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').Value = "'0.5162"
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').Value = "'0.3996"
$ws.Range('E8').Value = '  +1.71%  '
$ws.Range('D9').Value = "'0.08428"
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = "'42.68"
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('D11').Value = "'1.116"
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = "'23.33"
$ws.Range('E12').Value = '  +12.61%  '
$ws.Range('D13').Value = "'6.430"
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('D14').Value = '1.896.45'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = "'7.337"
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').Value = "'1.001"
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = "'94.17"
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('D19').Value = "'0.06640"
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('D20').Value = "'18.23"
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D21').Value = "'1.001"
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').Value = "'5.950"
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('D23').Value = '30.199.44'
$ws.Range('E23').Value = '  +2.92%  '
$ws.Range('D24').Value = "'11.28"
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('D25').Value = "'2.227"
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').Value = '2.110.29'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = "'21.68"
$ws.Range('E27').Value = '  +3.96%  '
$ws.Range('D29').Value = "'2.353"
$ws.Range('E29').Value = '  -3.16%  '
$ws.Range('D30').Value = "'129.05"
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('E31').Value = '  +3.27%  '
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').Value = "'6.106"
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').Value = "'3.757"
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('D35').Value = "'0.02501"
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('D36').Value = "'0.06544"
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('D37').Value = "'5.276"
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = "'1.218"
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('D40').Value = "'11.75"
$ws.Range('E40').Value = '  +4.33%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'8.728"
$ws.Range('E41').Value = '  -3.48%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = "'0.6495"
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').Value = "'1.231"
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').Value = "'0.6104"
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('D45').Value = "'13.22"
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').Value = "'3.702"
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('D47').Value = "'2.053"
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').Value = "'1.235"
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('D49').Value = "'124.51"
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('D50').Value = "'1.163"
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('D51').Value = "'78.98"
$ws.Range('E51').Value = '  +1.64%  '
